$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.558.50'
$ws.Range("E2").Value = '  +1.91%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.664.79'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.59'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4798'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2628'
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06156'
$ws.Range("E9").Value = '  +2.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07077'
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.666.28'
$ws.Range("E11").Value = '  +0.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.83'
$ws.Range("E12").Value = '  +2.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5925'
$ws.Range("E13").Value = '  -4.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.390'
$ws.Range("E14").Value = '  -4.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '74.44'
$ws.Range("E15").Value = '  +1.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.572.75'
$ws.Range("E18").Value = '  +1.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006760'
$ws.Range("E19").Value = '  +2.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.43'
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.879.39'
$ws.Range("E21").Value = '  +0.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.450'
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.686'
$ws.Range("E23").Value = '  +2.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.316'
$ws.Range("E24").Value = '  +1.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.61'
$ws.Range("E25").Value = '  +1.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.05'
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.406'
$ws.Range("E27").Value = '  +0.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '104.73'
$ws.Range("E28").Value = '  +3.24%  '
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.949'
$ws.Range("E30").Value = '  +4.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.662'
$ws.Range("E31").Value = '  +3.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07656'
$ws.Range("E32").Value = '  -3.25%  '
$ws.Range("B33").Value = 'Frax'
$ws.Range("C33").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9995'
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04325'
$ws.Range("E34").Value = '  -5.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.616'
$ws.Range("E35").Value = '  +0.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6137'
$ws.Range("E36").Value = '  +6.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9491'
$ws.Range("E37").Value = '  +1.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.613'
$ws.Range("E38").Value = '  -0.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8589'
$ws.Range("E39").Value = '  +2.06%  '
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01511'
$ws.Range("E41").Value = '  -1.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.880'
$ws.Range("E42").Value = '  +2.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.79'
$ws.Range("E43").Value = '  -0.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3766'
$ws.Range("E44").Value = '  +1.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.705'
$ws.Range("E45").Value = '  -2.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1120'
$ws.Range("E46").Value = '  +0.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.221'
$ws.Range("E47").Value = '  +2.72%  '
$ws.Range("E48").Value = '  +2.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.49'
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.373'
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("E51").Value = '  +0.15%  '
